# Auto-generated edit script: updates Leve-profit numeric columns (H-N)
# across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 780.4286
$ws.Range("I28").Value = 693
$ws.Range("K28").Value = 693
$ws.Range("M28").Value = -208
$ws.Range("H43").Value = 2704.4092
$ws.Range("I43").Value = 2874.8125
$ws.Range("J43").Value = 2250
$ws.Range("K43").Value = 2874.8125
$ws.Range("L43").Value = 2250
$ws.Range("M43").Value = -2805.8125
$ws.Range("N43").Value = -2388
$ws.Range("H51").Value = 3970886.2
$ws.Range("I51").Value = 55556556
$ws.Range("J51").Value = 2757.6924
$ws.Range("K51").Value = 55556556
$ws.Range("L51").Value = 2757.6924
$ws.Range("M51").Value = -55556072
$ws.Range("N51").Value = -3725.6924
$ws.Range("H53").Value = 276.13333
$ws.Range("I53").Value = 97.59999999999999
$ws.Range("J53").Value = 365.4
$ws.Range("K53").Value = 97.59999999999999
$ws.Range("L53").Value = 365.4
$ws.Range("M53").Value = 539.4
$ws.Range("N53").Value = -1639.4
$ws.Range("H76").Value = 6000.5454
$ws.Range("I76").Value = 3700.6
$ws.Range("J76").Value = 29000
$ws.Range("K76").Value = 3700.6
$ws.Range("L76").Value = 29000
$ws.Range("M76").Value = -3385.6
$ws.Range("N76").Value = -29630
$ws.Range("H79").Value = 6000.5454
$ws.Range("I79").Value = 3700.6
$ws.Range("J79").Value = 29000
$ws.Range("K79").Value = 3700.6
$ws.Range("L79").Value = 29000
$ws.Range("M79").Value = -2608.6
$ws.Range("N79").Value = -31184
$ws.Range("H86").Value = 3970.0527
$ws.Range("I86").Value = 3511.7273
$ws.Range("J86").Value = 4600.25
$ws.Range("K86").Value = 3511.7273
$ws.Range("L86").Value = 4600.25
$ws.Range("M86").Value = -2388.7273
$ws.Range("N86").Value = -6846.25
$ws.Range("H88").Value = 4699.6313
$ws.Range("I88").Value = 5768.846
$ws.Range("J88").Value = 2383
$ws.Range("K88").Value = 5768.846
$ws.Range("L88").Value = 2383
$ws.Range("M88").Value = -5362.846
$ws.Range("N88").Value = -3195
$ws.Range("H89").Value = 3970.0527
$ws.Range("I89").Value = 3511.7273
$ws.Range("J89").Value = 4600.25
$ws.Range("K89").Value = 17558.6365
$ws.Range("L89").Value = 23001.25
$ws.Range("M89").Value = -11942.6365
$ws.Range("N89").Value = -34233.25
$ws.Range("H91").Value = 4699.6313
$ws.Range("I91").Value = 5768.846
$ws.Range("J91").Value = 2383
$ws.Range("K91").Value = 5768.846
$ws.Range("L91").Value = 2383
$ws.Range("M91").Value = -4364.846
$ws.Range("N91").Value = -5191
$ws.Range("H107").Value = 860
$ws.Range("J107").Value = 810
$ws.Range("L107").Value = 810
$ws.Range("N107").Value = -4650
$ws.Range("H112").Value = 1184.1459
$ws.Range("J112").Value = 1198.674
$ws.Range("L112").Value = 3596.022
$ws.Range("N112").Value = -5812.022
$ws.Range("H121").Value = 1508.3334
$ws.Range("I121").Value = 347.5
$ws.Range("K121").Value = 1042.5
$ws.Range("M121").Value = 704.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28510.709
$ws.Range("I32").Value = 29807.365
$ws.Range("J32").Value = 21768.1
$ws.Range("K32").Value = 29807.365
$ws.Range("L32").Value = 21768.1
$ws.Range("M32").Value = -29520.365
$ws.Range("N32").Value = -22342.1
$ws.Range("H45").Value = 1434.3334
$ws.Range("I45").Value = 1378
$ws.Range("K45").Value = 1378
$ws.Range("M45").Value = -1001
$ws.Range("H132").Value = 11112604
$ws.Range("I132").Value = 16130066
$ws.Range("J132").Value = 2510.4285
$ws.Range("K132").Value = 48390198
$ws.Range("L132").Value = 7531.2855
$ws.Range("M132").Value = -48387668
$ws.Range("N132").Value = -12591.2855

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 204749.1
$ws.Range("J4").Value = 204749.1
$ws.Range("L4").Value = 204749.1
$ws.Range("N4").Value = -204973.1
$ws.Range("H31").Value = 167079.38
$ws.Range("I31").Value = 1907.909
$ws.Range("J31").Value = 223857.06
$ws.Range("K31").Value = 1907.909
$ws.Range("L31").Value = 223857.06
$ws.Range("M31").Value = -1612.909
$ws.Range("N31").Value = -224447.06
$ws.Range("H34").Value = 167079.38
$ws.Range("I34").Value = 1907.909
$ws.Range("J34").Value = 223857.06
$ws.Range("K34").Value = 1907.909
$ws.Range("L34").Value = 223857.06
$ws.Range("M34").Value = -1705.909
$ws.Range("N34").Value = -224261.06
$ws.Range("H100").Value = 25377
$ws.Range("J100").Value = 31754
$ws.Range("L100").Value = 31754
$ws.Range("N100").Value = -33918
$ws.Range("H111").Value = 28197.5
$ws.Range("J111").Value = 28197.5
$ws.Range("L111").Value = 28197.5
$ws.Range("N111").Value = -36377.5
$ws.Range("H115").Value = 34261
$ws.Range("J115").Value = 34261
$ws.Range("L115").Value = 34261
$ws.Range("N115").Value = -36611
$ws.Range("H125").Value = 23092.8
$ws.Range("J125").Value = 23092.8
$ws.Range("L125").Value = 23092.8
$ws.Range("N125").Value = -28012.8
$ws.Range("H132").Value = 56227.92
$ws.Range("I132").Value = 1885.1052
$ws.Range("J132").Value = 203729.86
$ws.Range("K132").Value = 5655.3156
$ws.Range("L132").Value = 611189.58
$ws.Range("M132").Value = -3125.3156
$ws.Range("N132").Value = -616249.58

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 934.6923
$ws.Range("J4").Value = 1650.1428
$ws.Range("L4").Value = 4950.428400000001
$ws.Range("N4").Value = -5174.428400000001
$ws.Range("H120").Value = 574180
$ws.Range("I120").Value = 601452
$ws.Range("J120").Value = 506000
$ws.Range("K120").Value = 1804356
$ws.Range("L120").Value = 1518000
$ws.Range("M120").Value = -1799518
$ws.Range("N120").Value = -1527676
$ws.Range("H122").Value = 6234.3687
$ws.Range("I122").Value = 545.5
$ws.Range("J122").Value = 10371.728
$ws.Range("K122").Value = 4909.5
$ws.Range("L122").Value = 93345.552
$ws.Range("M122").Value = -2459.5
$ws.Range("N122").Value = -98245.552
$ws.Range("H131").Value = 28297.338
$ws.Range("I131").Value = 8849.083000000001
$ws.Range("J131").Value = 31887.785
$ws.Range("K131").Value = 26547.249
$ws.Range("L131").Value = 95663.355
$ws.Range("M131").Value = -21507.249
$ws.Range("N131").Value = -105743.355
$ws.Range("H139").Value = 81840.87
$ws.Range("I139").Value = 252661.67
$ws.Range("J139").Value = 3000.5
$ws.Range("K139").Value = 757985.01
$ws.Range("L139").Value = 9001.5
$ws.Range("M139").Value = -752845.01
$ws.Range("N139").Value = -19281.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 16222.223
$ws.Range("J5").Value = 16375
$ws.Range("L5").Value = 16375
$ws.Range("N5").Value = -16599
$ws.Range("H95").Value = 22398.666
$ws.Range("J95").Value = 22398.666
$ws.Range("L95").Value = 22398.666
$ws.Range("N95").Value = -27890.666
$ws.Range("H110").Value = 32894.6
$ws.Range("J110").Value = 32894.6
$ws.Range("L110").Value = 32894.6
$ws.Range("N110").Value = -41074.6
$ws.Range("H120").Value = 39317
$ws.Range("J120").Value = 39317
$ws.Range("L120").Value = 39317
$ws.Range("N120").Value = -48993
$ws.Range("H126").Value = 5150.6
$ws.Range("I126").Value = 7379.1113
$ws.Range("K126").Value = 22137.3339
$ws.Range("M126").Value = -19667.3339
$ws.Range("H132").Value = 3205.25
$ws.Range("I132").Value = 1207.7858
$ws.Range("J132").Value = 7866
$ws.Range("K132").Value = 3623.3574
$ws.Range("L132").Value = 23598
$ws.Range("M132").Value = -1093.3574
$ws.Range("N132").Value = -28658
$ws.Range("H134").Value = 24311.111
$ws.Range("J134").Value = 24311.111
$ws.Range("L134").Value = 72933.333
$ws.Range("N134").Value = -78003.333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 19428.643
$ws.Range("H101").Value = 31976.834
$ws.Range("J101").Value = 31976.834
$ws.Range("L101").Value = 31976.834
$ws.Range("N101").Value = -38466.834
$ws.Range("H121").Value = 23952.666
$ws.Range("J121").Value = 23952.666
$ws.Range("L121").Value = 23952.666
$ws.Range("N121").Value = -27446.666
$ws.Range("H130").Value = 38830.25
$ws.Range("J130").Value = 38830.25
$ws.Range("L130").Value = 38830.25
$ws.Range("N130").Value = -48870.25
$ws.Range("H138").Value = 58899
$ws.Range("J138").Value = 58899
$ws.Range("L138").Value = 58899
$ws.Range("N138").Value = -69179

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3375.5
$ws.Range("H56").Value = 11271.333
$ws.Range("J56").Value = 11271.333
$ws.Range("L56").Value = 11271.333
$ws.Range("N56").Value = -12699.333
$ws.Range("H76").Value = 27696.666
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 27696.666
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 27696.666
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -28326.666
$ws.Range("H79").Value = 27696.666
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 27696.666
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 27696.666
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -29880.666
$ws.Range("H95").Value = 40773.332
$ws.Range("J95").Value = 40773.332
$ws.Range("L95").Value = 40773.332
$ws.Range("N95").Value = -46265.332
$ws.Range("H103").Value = 35261.332
$ws.Range("J103").Value = 35261.332
$ws.Range("L103").Value = 35261.332
$ws.Range("N103").Value = -37605.332
$ws.Range("H113").Value = 1183.3334
$ws.Range("I113").Value = 1183.3334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3550.0002
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1380.0002
$ws.Range("N113").ClearContents()
$ws.Range("H121").Value = 44416
$ws.Range("J121").Value = 44416
$ws.Range("L121").Value = 44416
$ws.Range("N121").Value = -47910
$ws.Range("H125").Value = 39715
$ws.Range("J125").Value = 39715
$ws.Range("L125").Value = 39715
$ws.Range("N125").Value = -49555
$ws.Range("H138").Value = 45357.25
$ws.Range("J138").Value = 45357.25
$ws.Range("L138").Value = 45357.25
$ws.Range("N138").Value = -55637.25

